$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-14 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-15 Monday", 2)

$d.Content.Find.Execute("950÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "851÷3=", 2)
$d.Content.Find.Execute("356÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "297÷2=", 2)
$d.Content.Find.Execute("785÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "153÷5=", 2)
$d.Content.Find.Execute("558÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "795÷3=", 2)
$d.Content.Find.Execute("187÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷6=", 2)
$d.Content.Find.Execute("663÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "474÷8=", 2)
$d.Content.Find.Execute("486÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "944÷8=", 2)
$d.Content.Find.Execute("964÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "843÷8=", 2)
$d.Content.Find.Execute("102÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "288÷7=", 2)
$d.Content.Find.Execute("128÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "331÷9=", 2)
$d.Content.Find.Execute("854÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "966÷8=", 2)
$d.Content.Find.Execute("857÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "785÷6=", 2)
$d.Content.Find.Execute("614÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "851÷2=", 2)
$d.Content.Find.Execute("520÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "781÷4=", 2)
$d.Content.Find.Execute("990÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "192÷6=", 2)
$d.Content.Find.Execute("553÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "867÷2=", 2)
$d.Content.Find.Execute("823÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "885÷3=", 2)
$d.Content.Find.Execute("663÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "133÷6=", 2)
$d.Content.Find.Execute("758÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "836÷3=", 2)
$d.Content.Find.Execute("153÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "357÷4=", 2)
$d.Content.Find.Execute("825÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "166÷5=", 2)
$d.Content.Find.Execute("330÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "945÷2=", 2)
$d.Content.Find.Execute("307÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "506÷8=", 2)
$d.Content.Find.Execute("304÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "949÷2=", 2)
$d.Content.Find.Execute("529÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "103÷2=", 2)
